# The site footer ("Ver no Jupiter Salvar em pdf Salvar em docx" and the
# "© 2020 ..." copyright line), together with the blank paragraph that
# separates them from the page break above, is removed from the end of
# the document. The paragraph carrying the page break itself, and the
# blank paragraph that used to follow the copyright line, are left in
# place (and end up adjacent to each other).

$d = $word.ActiveDocument

# Locate the "Ver no Jupiter ..." paragraph - the start of the block to remove.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Ver no Jupiter*") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    # One blank paragraph precedes it, the copyright paragraph follows it -
    # both are part of the block being deleted.
    $blockStart = $target.Previous(1)
    $blockEnd = $target.Next(1)

    $range = $d.Range($blockStart.Range.Start, $blockEnd.Range.End)
    $range.Delete()
}
